$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Post-Test (column C) values for rows 2-13
$values = @{
    2  = 60
    3  = 51
    4  = 54
    5  = 37
    6  = 52
    7  = 40
    8  = 53
    9  = 52
    10 = 54
    11 = 45
    12 = 48
    13 = 48
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

# Update the active selection to C14 (matching the recorded sheet view state)
$ws.Range("C14").Select()
